$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(100).OutlineLevel = 3
Write-Output ("Row100 outlinelevel=" + $ws.Rows(100).OutlineLevel())
$ws.Rows(46).OutlineLevel = 2
Write-Output ("Row46 outlinelevel=" + $ws.Rows(46).OutlineLevel())
